$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 previously had empty A3/B3 (style already set) and no C3 cell at all.
# Set the new values for A3 and B3 (they keep their existing style).
$ws.Range("A3").Value = "viji@123"
$ws.Range("B3").Value = "viji@123"

# C3 is a brand-new cell; copy the formatting from B3 (style index 2)
# before assigning its value so it matches A3/B3's style exactly.
$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C3").Value = "vijayarani"

$excel.CutCopyMode = 0
